$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B3").Value = 2772.04168703704
$ws.Range("D8").Value = 724.6697314911827
